# Apply "Practice tasks and final revisions" edits to the task-order workbook.
# Renames each sheet (new timestamp suffixes) and updates the stimulus/order
# file names referenced in column B of every sheet.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets -----------------------------------------------------
$sheetRenames = @(
    @{ Index = 1; Old = "GNG_TO-16504778363739493";  New = "GNG_TO-16509961261482885"  },
    @{ Index = 2; Old = "NB_TO-16504778385419507";   New = "NB_TO-16509961290847852"   },
    @{ Index = 3; Old = "RS_TO-16504778385479517";   New = "RS_TO-16509961290847852"   },
    @{ Index = 4; Old = "TOL_TO-16504778386039479";  New = "TOL_TO-16509961291487718"  },
    @{ Index = 5; Old = "vSAT_TO-16504778386669817"; New = "vSAT_TO-16509961292127697" }
)

foreach ($rename in $sheetRenames) {
    $sheet = $wb.Worksheets.Item($rename.Index)
    $sheet.Name = $rename.New
}

# --- Update cell values on sheet 1 (GNG_TO) --------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961261082761.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961261322777.csv"
$ws1.Range("B4").Value = "go_stims-16509961261322777.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961261482885.csv"

# --- Update cell values on sheet 2 (NB_TO) ---------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-165099612815644.csv"
$ws2.Range("B3").Value = "TB-16509961289004436.csv"
$ws2.Range("B4").Value = "OB-16509961279244144.csv"
$ws2.Range("B5").Value = "TB-16509961290687296.csv"
$ws2.Range("B6").Value = "ZB-match_7-165099612637228.csv"
$ws2.Range("B7").Value = "ZB-match_5-1650996127092271.csv"
$ws2.Range("B8").Value = "OB-165099612773244.csv"
$ws2.Range("B9").Value = "ZB-match_4-16509961269642725.csv"
$ws2.Range("B10").Value = "OB-165099612718831.csv"

# --- Sheet 3 (RS_TO) has no data changes, only the sheet name above -------

# --- Update cell values on sheet 4 (TOL_TO) --------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961291167567.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961290927253.csv"
$ws4.Range("B4").Value = "MM_stims-16509961291327245.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961291167567.csv"
$ws4.Range("B6").Value = "MM_stims-16509961291487718.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961291327245.csv"

# --- Update cell values on sheet 5 (vSAT_TO) -------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961291807292.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961291967747.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961291647651.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961291487718.csv"
